$wb = $excel.ActiveWorkbook

# --- Sheet "Aprilie": update B4, C4, D4 with corrected results ---
$wsAprilie = $wb.Worksheets.Item("Aprilie")
$wsAprilie.Range("B4").Value = 96.69
$wsAprilie.Range("C4").Value = 96.6
$wsAprilie.Range("D4").Value = 96.7

# --- Sheet "16-mai": update B4, C4, D4 with corrected results, and move selection ---
$ws16mai = $wb.Worksheets.Item("16-mai")
$ws16mai.Range("B4").Value = 97.46
$ws16mai.Range("C4").Value = 97.3
$ws16mai.Range("D4").Value = 97.5

$ws16mai.Activate() | Out-Null
$ws16mai.Range("E10").Select() | Out-Null
